$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; this shifts the existing "District" column
# (F) to the right, into G, and leaves a blank column F behind.
$ws.Columns.Item(6).Insert()

# New header for the inserted column.
$ws.Cells.Item(2, 6).Value2 = "Address"

# Find the last populated data row (column A holds the running serial
# number for every data row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $nameAndAddress = $ws.Cells.Item($r, 2).Value2

    $address = ""

    if ($nameAndAddress) {
        $lines = $nameAndAddress -split "`n"

        if ($lines.Length -eq 2) {
            $addrLine = $lines[1].TrimEnd(".")
            $parts = $addrLine -split ","

            if ($parts.Length -ge 2) {
                $trimmedParts = @()
                for ($i = 0; $i -lt $parts.Length - 1; $i++) {
                    $trimmedParts += $parts[$i].Trim()
                }
                $address = $trimmedParts -join ""
            }
        }
    }

    if ($address -ne "") {
        $ws.Cells.Item($r, 6).Value2 = $address
    }
}
